$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.679.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.01%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.564.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.76%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '652.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.81%  '

# Row 7
$ws.Range("E7").Value = '  +1.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.402'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.49%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.995'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.21%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.563.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.85%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.202'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.43%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.19%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.246.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.11%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.469.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.00%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000253'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.71%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.572.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.53%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.53%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '506.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.72%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.478'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.40%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000194'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.46%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.44%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.10%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.752.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.04%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.12%  '

# Row 31
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.73%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.139'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.40%  '

# Row 34
$ws.Range("E34").Value = '  +0.09%  '

# Row 35
$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.177'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.40%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.560'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.55%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '559.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.80%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.150'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.54%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.903'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '35.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +39.81%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.63%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.17%  '

# Row 47
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.72%  '

# Row 48
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.29%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0412'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.31%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.27%  '
